# Revamp PricingDownload: append new pricing line items (rows 11-18)
# to the existing order sheet. Columns: A=SKU, B=Item, C=Quantity,
# D=Cost Per, E=Total Cost. All values are stored as text in the
# original workbook (t="inlineStr"/shared-string), so numeric-looking
# columns (C/D/E) are entered with a leading apostrophe to force text
# entry, then the cell style is reset to "Normal" so no stray
# quote-prefix style sticks to the cell (matching the source file,
# which carries no style attribute on these cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 11; Sku = "SLOPL4";       Item = "Lid - Portion (3.25oz)";              Qty = "1"; Cost = "25.90";  Total = "25.90"  },
    @{ Row = 12; Sku = "TS12";         Item = "Tamper Evident - 12oz Square";        Qty = "1"; Cost = "38.39";  Total = "38.39"  },
    @{ Row = 13; Sku = "TS16";         Item = "Tamper Evident - 16oz";               Qty = "1"; Cost = "41.87";  Total = "41.87"  },
    @{ Row = 14; Sku = "TS8";          Item = "Tamper Evident - 8oz";                Qty = "2"; Cost = "38.30";  Total = "76.60"  },
    @{ Row = 15; Sku = "ANPM424";      Item = "Container - Anchor (24oz)";           Qty = "2"; Cost = "47.17";  Total = "94.34"  },
    @{ Row = 16; Sku = "6G063015";     Item = "Bag Poly - 6x3x15 LW";                Qty = "2"; Cost = "15.06";  Total = "30.12"  },
    @{ Row = 17; Sku = "HIMF1824XC";   Item = "Bag Poly - 10x8x24 HW (Soup Bag)";    Qty = "1"; Cost = "63.33";  Total = "63.33"  },
    @{ Row = 18; Sku = "SAB52032T300"; Item = "Lid Salad - 24/32oz Sabert (Round)";  Qty = "2"; Cost = "80.81";  Total = "161.62" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.Sku
    $ws.Range("B$row").Value = $r.Item

    # Quantity / Cost Per / Total Cost look numeric, so force text entry
    # with a quote prefix, then strip the resulting quote-prefix style
    # back to Normal so the cell ends up with plain text + default style.
    $ws.Range("C$row").Value = "'" + $r.Qty
    $ws.Range("C$row").Style = "Normal"

    $ws.Range("D$row").Value = "'" + $r.Cost
    $ws.Range("D$row").Style = "Normal"

    $ws.Range("E$row").Value = "'" + $r.Total
    $ws.Range("E$row").Style = "Normal"
}
